$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.657.64"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.204.41"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.22"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.98"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -4.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.12"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "2.534.82"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "2.210.75"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "42.543.74"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.95"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.18"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.99"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.37"
$ws.Range("E24").Value = "  -8.38%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.52"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.09"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.11"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0837"
$ws.Range("E33").Value = "  +5.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.42"
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  +18.37%  "
$ws.Range("E42").Value = "  -6.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.70"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.35"
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.460"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("E51").Value = "  -1.12%  "
